# AutoCommit_10 ноября 2023 г. 11:27:34_SibNout2023
#
# Marks additional "ок" (checkmark/credit) grades for students in rows 7,
# 26 and 27 of the gradebook sheet, matching the homework columns that were
# already filled in for the other students.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mark = "ок"

# --- Row 7 (Боярский Артём): дз1, д4, дз2, д3 and д5, д6 are done; ок column (G) stays blank ---
$ws.Range("C7").Value = $mark
$ws.Range("D7").Value = $mark
$ws.Range("E7").Value = $mark
$ws.Range("F7").Value = $mark
$ws.Range("H7").Value = $mark

# I7 did not exist yet in this row - create it with the same value, then
# pull the formatting (borders/alignment/font) from the neighbouring cell
# so it lands on the same shared cell style instead of the default one.
$ws.Range("I7").Value = $mark
$ws.Range("H7").Copy()
$ws.Range("I7").PasteSpecial(-4122)

# --- Row 26 (Писецкий Михаил): every homework column is now complete ---
$ws.Range("C26").Value = $mark
$ws.Range("D26").Value = $mark
$ws.Range("E26").Value = $mark
$ws.Range("F26").Value = $mark
$ws.Range("G26").Value = $mark
$ws.Range("H26").Value = $mark

$ws.Range("I26").Value = $mark
$ws.Range("H26").Copy()
$ws.Range("I26").PasteSpecial(-4122)

# --- Row 27 (Подлесный Никита): all but column F are done ---
$ws.Range("C27").Value = $mark
$ws.Range("D27").Value = $mark
$ws.Range("E27").Value = $mark
$ws.Range("G27").Value = $mark
$ws.Range("H27").Value = $mark

# Leave the cursor on the last cell that was edited, like the original author.
$ws.Range("I7").Select()
